# Weekly fruit/vegetable price update: insert a new week's worth of
# "Repollo" (cabbage) price records for Vega Central Mapocho de Santiago
# just above the existing row 662, shifting the rest of the table down.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 blank rows at 662..665 (existing rows 662+ shift down to 666+)
$ws.Range("A662:A665").EntireRow.Insert()

# Row 662: Crespo record / Primera
$ws.Range("A662").Value = 9
$ws.Range("B662").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C662").Value = "Metropolitana"
$ws.Range("D662").Value = 45013
$ws.Range("E662").Value = 13
$ws.Range("F662").Value = 100112006
$ws.Range("G662").Value = "Repollo"
$ws.Range("H662").Value = "Crespo record"
$ws.Range("I662").Value = "Primera"
$ws.Range("J662").Value = 3400
$ws.Range("K662").Value = 1300
$ws.Range("L662").Value = 1400
$ws.Range("M662").Value = 1350
$ws.Range("N662").Value = "`$/unidad"
$ws.Range("O662").Value = "Provincia de Quillota"
$ws.Range("P662").Value = 1350
$ws.Range("Q662").Value = 1
$ws.Range("R662").Value = "Hortaliza"

# Row 663: Crespo record / Segunda
$ws.Range("A663").Value = 9
$ws.Range("B663").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C663").Value = "Metropolitana"
$ws.Range("D663").Value = 45013
$ws.Range("E663").Value = 13
$ws.Range("F663").Value = 100112006
$ws.Range("G663").Value = "Repollo"
$ws.Range("H663").Value = "Crespo record"
$ws.Range("I663").Value = "Segunda"
$ws.Range("J663").Value = 1600
$ws.Range("K663").Value = 1000
$ws.Range("L663").Value = 1000
$ws.Range("M663").Value = 1000
$ws.Range("N663").Value = "`$/unidad"
$ws.Range("O663").Value = "Provincia de Quillota"
$ws.Range("P663").Value = 1000
$ws.Range("Q663").Value = 1
$ws.Range("R663").Value = "Hortaliza"

# Row 664: Morada(o) / Primera
$ws.Range("A664").Value = 9
$ws.Range("B664").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C664").Value = "Metropolitana"
$ws.Range("D664").Value = 45013
$ws.Range("E664").Value = 13
$ws.Range("F664").Value = 100112006
$ws.Range("G664").Value = "Repollo"
$ws.Range("H664").Value = "Morada(o)"
$ws.Range("I664").Value = "Primera"
$ws.Range("J664").Value = 1600
$ws.Range("K664").Value = 1400
$ws.Range("L664").Value = 1500
$ws.Range("M664").Value = 1450
$ws.Range("N664").Value = "`$/unidad"
$ws.Range("O664").Value = "Provincia de Quillota"
$ws.Range("P664").Value = 1450
$ws.Range("Q664").Value = 1
$ws.Range("R664").Value = "Hortaliza"

# Row 665: Morada(o) / Segunda
$ws.Range("A665").Value = 9
$ws.Range("B665").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C665").Value = "Metropolitana"
$ws.Range("D665").Value = 45013
$ws.Range("E665").Value = 13
$ws.Range("F665").Value = 100112006
$ws.Range("G665").Value = "Repollo"
$ws.Range("H665").Value = "Morada(o)"
$ws.Range("I665").Value = "Segunda"
$ws.Range("J665").Value = 790
$ws.Range("K665").Value = 1200
$ws.Range("L665").Value = 1200
$ws.Range("M665").Value = 1200
$ws.Range("N665").Value = "`$/unidad"
$ws.Range("O665").Value = "Provincia de Quillota"
$ws.Range("P665").Value = 1200
$ws.Range("Q665").Value = 1
$ws.Range("R665").Value = "Hortaliza"

# Keep the D column's date-style number format consistent with the rest
# of the column (style carries over from the insert, but set explicitly
# to be safe).
$ws.Range("D662:D665").NumberFormat = $ws.Range("D661").NumberFormat
